# GroupePR.xlsx update
# - Removes the obsolete "2DEVWFS D2" group row.
# - Adds a new "groupe_physique" column (F) mirroring the "Group Présentiel"
#   column (A) for every remaining data row.
# - Updates the sheet selection to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "2DEVWFS D2" row (row 4); remaining rows shift up.
$ws.Rows(4).Delete()

# New header for column F.
$ws.Cells.Item(1, 6).Value = "groupe_physique"

# Mirror column A into the new column F for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 1).Value2
}

# Match the print orientation recorded by Excel on save.
$ws.PageSetup.Orientation = 1

# Reflect the final selection/scroll state left after the edit.
$ws.Range("F27:F35").Select() | Out-Null

Write-Host "GroupePR.xlsx updated: removed '2DEVWFS D2' row, added groupe_physique column"
